$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells remain text (matching original inlineStr formatting),
# so numeric-looking strings (e.g. "54.50", "9.00") are not coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.897.85'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +5.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.533.80'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +5.67%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +6.10%  '
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '188.14'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +9.31%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.525.49'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +5.63%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.632'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.10%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +14.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.50'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000270'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.33'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.099.90'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.532.67'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +5.31%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.52'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +5.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '66.883.06'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.07'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.89%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '423.92'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +13.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.13'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +10.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.51'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.19'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.08'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.91'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +8.06%  '
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.28'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +8.97%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.10'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.00'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +9.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.41'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +5.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '630.78'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.67'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.74'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.07%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '60.08'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '38.28'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0812'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +12.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.148'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +18.70%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.128.67'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.31'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +9.78%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.80%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.88'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +10.60%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.37'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +9.37%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.28%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.06%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '140.97'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.71%  '
